$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 486
$ws.Range("I12").Value = 419.5
$ws.Range("K12").Value = 419.5
$ws.Range("M12").Value = -249.5
$ws.Range("H19").Value = 1657.6154
$ws.Range("I19").Value = 1695.75
$ws.Range("K19").Value = 1695.75
$ws.Range("M19").Value = -1520.75
$ws.Range("H53").Value = 1136.6666
$ws.Range("I53").Value = 647.3333
$ws.Range("K53").Value = 647.3333
$ws.Range("M53").Value = -10.33330000000001
$ws.Range("H107").Value = 1251.7812
$ws.Range("I107").Value = 1380.375
$ws.Range("J107").Value = 866
$ws.Range("K107").Value = 1380.375
$ws.Range("L107").Value = 866
$ws.Range("M107").Value = 539.625
$ws.Range("N107").Value = -4706
$ws.Range("H112").Value = 103690.6
$ws.Range("J112").Value = 114422.89
$ws.Range("L112").Value = 343268.67
$ws.Range("N112").Value = -345484.67
$ws.Range("H131").Value = 3203.6924
$ws.Range("I131").Value = 1264.2858
$ws.Range("J131").Value = 5466.3335
$ws.Range("K131").Value = 3792.8574
$ws.Range("L131").Value = 16399.0005
$ws.Range("M131").Value = 1247.1426
$ws.Range("N131").Value = -26479.0005
$ws.Range("H132").Value = 1400.7188
$ws.Range("I132").Value = 1462.5
$ws.Range("J132").Value = 803.5
$ws.Range("K132").Value = 4387.5
$ws.Range("L132").Value = 2410.5
$ws.Range("M132").Value = -1857.5
$ws.Range("N132").Value = -7470.5
$ws.Range("H137").Value = 41989.383
$ws.Range("I137").Value = 84468.89999999999
$ws.Range("K137").Value = 253406.7
$ws.Range("M137").Value = -250856.7
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 93
$ws.Range("I4").Value = 99
$ws.Range("K4").Value = 99
$ws.Range("M4").Value = 17
$ws.Range("H32").Value = 23824520
$ws.Range("I32").Value = 30130784
$ws.Range("K32").Value = 30130784
$ws.Range("M32").Value = -30130497
$ws.Range("H63").Value = 4079.4443
$ws.Range("I63").Value = 2428.125
$ws.Range("J63").Value = 5400.5
$ws.Range("K63").Value = 2428.125
$ws.Range("L63").Value = 5400.5
$ws.Range("M63").Value = -1742.125
$ws.Range("N63").Value = -6772.5
$ws.Range("H66").Value = 4079.4443
$ws.Range("I66").Value = 2428.125
$ws.Range("J66").Value = 5400.5
$ws.Range("K66").Value = 12140.625
$ws.Range("L66").Value = 27002.5
$ws.Range("M66").Value = -8708.625
$ws.Range("N66").Value = -33866.5
$ws.Range("H74").Value = 2421.359
$ws.Range("I74").Value = 2207.1765
$ws.Range("J74").Value = 3877.8
$ws.Range("K74").Value = 2207.1765
$ws.Range("L74").Value = 3877.8
$ws.Range("M74").Value = -1333.1765
$ws.Range("N74").Value = -5625.8
$ws.Range("H77").Value = 2421.359
$ws.Range("I77").Value = 2207.1765
$ws.Range("J77").Value = 3877.8
$ws.Range("K77").Value = 11035.8825
$ws.Range("L77").Value = 19389
$ws.Range("M77").Value = -6667.8825
$ws.Range("N77").Value = -28125
$ws.Range("H102").Value = 1603.4706
$ws.Range("I102").Value = 1125.7142
$ws.Range("K102").Value = 1125.7142
$ws.Range("M102").Value = 496.2858000000001
$ws.Range("H126").Value = 9997
$ws.Range("I126").Value = 9997
$ws.Range("K126").Value = 29991
$ws.Range("M126").Value = -27521
$ws.Range("H132").Value = 2850.2888
$ws.Range("I132").Value = 2415.0938
$ws.Range("K132").Value = 7245.2814
$ws.Range("M132").Value = -4715.2814
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 79330
$ws.Range("J116").Value = 79330
$ws.Range("L116").Value = 79330
$ws.Range("N116").Value = -88508
$ws.Range("H134").Value = 2860777.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 45996.668
$ws.Range("J68").Value = 74990
$ws.Range("L68").Value = 74990
$ws.Range("N68").Value = -76488
$ws.Range("H71").Value = 45996.668
$ws.Range("J71").Value = 74990
$ws.Range("L71").Value = 224970
$ws.Range("N71").Value = -232458
$ws.Range("H94").Value = 1328.2106
$ws.Range("I94").Value = 288.2
$ws.Range("K94").Value = 288.2
$ws.Range("M94").Value = 162.8
$ws.Range("H118").Value = 112494.5
$ws.Range("J118").Value = 105989
$ws.Range("L118").Value = 105989
$ws.Range("N118").Value = -109303
$ws.Range("H132").Value = 3530.7812
$ws.Range("I132").Value = 3338.348
$ws.Range("J132").Value = 4022.5557
$ws.Range("K132").Value = 10015.044
$ws.Range("L132").Value = 12067.6671
$ws.Range("M132").Value = -7485.044
$ws.Range("N132").Value = -17127.6671
$ws.Range("H134").Value = 1871.3125
$ws.Range("I134").Value = 1751.6333
$ws.Range("K134").Value = 5254.8999
$ws.Range("M134").Value = -2719.8999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 5494.25
$ws.Range("I18").Value = 4988.75
$ws.Range("K18").Value = 14966.25
$ws.Range("M18").Value = -14797.25
$ws.Range("H109").Value = 4150
$ws.Range("I109").Value = 3950
$ws.Range("J109").Value = 4250
$ws.Range("K109").Value = 11850
$ws.Range("L109").Value = 12750
$ws.Range("N109").Value = -14830
$ws.Range("M109").Value = -10810
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 106989
$ws.Range("J116").Value = 106989
$ws.Range("L116").Value = 106989
$ws.Range("N116").Value = -116167
$ws.Range("H124").Value = 153199
$ws.Range("J124").Value = 153199
$ws.Range("L124").Value = 153199
$ws.Range("N124").Value = -163019
$ws.Range("H126").Value = 3014.5
$ws.Range("I126").Value = 2917
$ws.Range("K126").Value = 8751
$ws.Range("M126").Value = -6281
$ws.Range("H132").Value = 3103.468
$ws.Range("I132").Value = 2980.2778
$ws.Range("K132").Value = 8940.8334
$ws.Range("M132").Value = -6410.8334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1419.3
$ws.Range("I61").Value = 1449.25
$ws.Range("K61").Value = 1449.25
$ws.Range("M61").Value = -1247.25
$ws.Range("H93").Value = 2236.75
$ws.Range("I93").Value = 1278.8
$ws.Range("K93").Value = 1278.8
$ws.Range("M93").Value = -30.79999999999995
$ws.Range("H113").Value = 1419.3
$ws.Range("I113").Value = 1449.25
$ws.Range("K113").Value = 1449.25
$ws.Range("M113").Value = 720.75
$ws.Range("H122").Value = 14312.833
$ws.Range("I122").Value = 14638.889
$ws.Range("K122").Value = 43916.667
$ws.Range("M122").Value = -41466.667
$ws.Range("H132").Value = 116199.11
$ws.Range("I132").Value = 171151.83
$ws.Range("J132").Value = 6293.6665
$ws.Range("K132").Value = 513455.49
$ws.Range("L132").Value = 18880.9995
$ws.Range("M132").Value = -510925.49
$ws.Range("N132").Value = -23940.9995
$ws.Range("H136").Value = 3404.6428
$ws.Range("I136").Value = 3226.4546
$ws.Range("K136").Value = 9679.363799999999
$ws.Range("M136").Value = -7129.363799999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 103678.75
$ws.Range("J16").Value = 103678.75
$ws.Range("L16").Value = 103678.75
$ws.Range("N16").Value = -104262.75
$ws.Range("H62").Value = 3949.6365
$ws.Range("J62").Value = 4212
$ws.Range("L62").Value = 4212
$ws.Range("N62").Value = -5460
$ws.Range("H64").Value = 99984
$ws.Range("J64").Value = 99984
$ws.Range("L64").Value = 99984
$ws.Range("N64").Value = -100480
$ws.Range("H65").Value = 3949.6365
$ws.Range("J65").Value = 4212
$ws.Range("L65").Value = 21060
$ws.Range("N65").Value = -27300
$ws.Range("H67").Value = 99984
$ws.Range("J67").Value = 99984
$ws.Range("L67").Value = 99984
$ws.Range("N67").Value = -101700
$ws.Range("H107").Value = 588.63635
$ws.Range("I107").Value = 467.7143
$ws.Range("J107").Value = 800.25
$ws.Range("K107").Value = 1403.1429
$ws.Range("L107").Value = 2400.75
$ws.Range("M107").Value = 516.8571000000002
$ws.Range("N107").Value = -6240.75
$ws.Range("H132").Value = 3805.4375
$ws.Range("I132").Value = 3646
$ws.Range("K132").Value = 10938
$ws.Range("M132").Value = -8408
